$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9 data
$ws.Range("A9").Value = "janluy moreno"
$ws.Range("B9").Value = "Estudiante"
$ws.Range("C9").Value = "janluy_moreno@cun.edu.co"
$ws.Range("D9").Value = "CC"
$ws.Range("E9").Value = 1022348425
$ws.Range("F9").Value = "Curso de ortografia"

# G9 must stay text (same as the rest of the column) rather than being
# auto-parsed into a date serial number.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "2023-06-07"

$ws.Range("H9").Value = 45
$ws.Range("I9").Value = "virtual"
$ws.Range("J9").Value = "Bogota"

# C9 becomes a mailto hyperlink (Excel auto-creates the built-in
# "Hyperlink" cell style + font the first time a hyperlink is inserted).
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:janluy_moreno@cun.edu.co", "", "", "janluy_moreno@cun.edu.co")

$ws.Range("H9:J9").Select()
